$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.660.08"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.277.75"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "123.91"
$ws.Range("E5").Value = "  +6.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.48"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.17"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.36"
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.44"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.905"
$ws.Range("E15").Value = "  +3.74%  "
$ws.Range("D16").Value = "2.623.18"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "2.274.79"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "43.682.21"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.96"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.31"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.42"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.03"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.49"
$ws.Range("E25").Value = "  -6.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.97"
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.19"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.69"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.65"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0918"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.37"
$ws.Range("E35").Value = "  +13.11%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("E37").Value = "  +4.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.93"
$ws.Range("E41").Value = "  -4.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.79"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.238"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.64"
$ws.Range("E46").Value = "  -11.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.11"
$ws.Range("E47").Value = "  +37.80%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.57"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.85"
$ws.Range("E51").Value = "  -0.94%  "
